$wb = $excel.ActiveWorkbook

# Remember which sheet was active before the edit so the workbook's
# selected-tab view stays unchanged (the diff only touches the <sheets>
# list, not <bookViews>).
$origActiveName = $wb.ActiveSheet.Name

# Add a new worksheet "mateusz" after the last existing sheet (emre),
# matching Mateusz's offer/quote data (see commit message).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "mateusz"

$data = @(
    @("groupComponent", "nameComponent"),
    @("RAM", "D416GB"),
    @("GPU", "GTX Titan X"),
    @("PSU", "Xpredator 750M"),
    @("Drive", "950 EVO"),
    @("Motherboard", "X99 Rampage V Extreme"),
    @("CPU", "i7 5960X")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}

$wb.Worksheets.Item($origActiveName).Activate()
